# bulk orders updated. Place Customer orders. Checkout.
# Rebuild the bulkOrders sheet with the new ItemCode/ItemName/... columns,
# two product rows (Leva S / Linea Mini), the refreshed header styling and
# the widened column layout.

function RGBColor($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values (row-major order so new shared strings line up with the
#    order headers/data are encountered on the sheet).
# ---------------------------------------------------------------------

# Row 1 - headers
$ws.Range("A1").Value = "ItemCode"
$ws.Range("B1").Value = "ItemName"
$ws.Range("C1").Value = "Itemtype"
$ws.Range("D1").Value = "FamilyCode"
$ws.Range("E1").Value = "FamilyDesc"
$ws.Range("F1").Value = "Isactive"
$ws.Range("G1").Value = "Weight"
$ws.Range("H1").Value = "WeightUoM"
$ws.Range("I1").Value = "Itemgroup"
$ws.Range("J1").Value = "GroupDescription"
$ws.Range("K1").Value = "ModelNumber "
$ws.Range("L1").Value = "Quantity"

# Row 2 - Leva S
$ws.Range("A2").Value = "levas123"
$ws.Range("B2").Value = "Leva S"
$ws.Range("C2").Value = "Machine"
$ws.Range("D2").Value = "Leva S"
$ws.Range("E2").Value = "Commercial Machine"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = 5000
$ws.Range("H2").Value = "Grams"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "levas123"
$ws.Range("L2").Value = ""

# Row 3 - Linea Mini
$ws.Range("A3").Value = "lineamini123"
$ws.Range("B3").Value = "Linea Mini"
$ws.Range("C3").Value = "Machine"
$ws.Range("D3").Value = "Linea Mini"
$ws.Range("E3").Value = "Commercial Machine"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = 4500
$ws.Range("H3").Value = "Grams"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "lineamini123"
$ws.Range("L3").Value = 3

# ---------------------------------------------------------------------
# 2. Restyle the A1 header template: bold Times New Roman, centered +
#    wrapped text, refreshed fill colors. Done on a single cell first so
#    the whole header row ends up sharing one consistent style.
# ---------------------------------------------------------------------

$tpl = $ws.Range("A1")
$tpl.Font.Bold = $true
$tpl.Font.Name = "Times New Roman"
$tpl.Font.Family = 1
$tpl.HorizontalAlignment = -4108
$tpl.WrapText = $true
$tpl.Interior.Color = RGBColor 0xEE 0xEE 0xEE
$tpl.Interior.PatternColor = RGBColor 0xFF 0xFF 0xCC

# ---------------------------------------------------------------------
# 3. Propagate formatting: the finished header style across B1:L1, and
#    the existing data-row style across the newly used columns C:L.
# ---------------------------------------------------------------------

$tpl.Copy()
$ws.Range("B1:L1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("C2:L2").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("C3:L3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Column widths for the new layout.
# ---------------------------------------------------------------------

$ws.Columns("A:D").ColumnWidth = 10.51
$ws.Columns("E:E").ColumnWidth = 17.52
$ws.Columns("F:I").ColumnWidth = 10.51
$ws.Columns("J:J").ColumnWidth = 13.61
$ws.Columns("K:K").ColumnWidth = 12.13
$ws.Columns("L:L").ColumnWidth = 13.21

# ---------------------------------------------------------------------
# 5. Selection matches the post-edit active cell.
# ---------------------------------------------------------------------

[void]$ws.Range("L2").Select()
